$d = $word.ActiveDocument
$p = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt spacing from text
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5

# Increase the left indent from 6pt (120 twips) to 11.25pt (225 twips)
$p.Format.LeftIndent = 11.25

# Update the placeholder id and drop the trailing space run
$d.Content.Find.Execute("**ID__AFFARS_5333_topic_6__ID** ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "**ID__AFFARS_5333_105__ID**", 2)
